$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 20
$ws.Range("H20").Value = 6437.25
$ws.Range("I20").Value = 6437.25
$ws.Range("K20").Value = 6437.25
$ws.Range("M20").Value = -6207.25
# row 35
$ws.Range("H35").Value = 6437.25
$ws.Range("I35").Value = 6437.25
$ws.Range("K35").Value = 6437.25
$ws.Range("M35").Value = -6058.25
# row 45
$ws.Range("H45").Value = 2682.6
$ws.Range("I45").Value = 1237
$ws.Range("J45").Value = 3646.3333
$ws.Range("K45").Value = 3711
$ws.Range("L45").Value = 10938.9999
$ws.Range("M45").Value = -3519
$ws.Range("N45").Value = -11322.9999
# row 70
$ws.Range("H70").Value = 12567524
$ws.Range("I70").Value = 146199.58
$ws.Range("K70").Value = 438598.74
$ws.Range("M70").Value = -438328.74
# row 73
$ws.Range("H73").Value = 12567524
$ws.Range("I73").Value = 146199.58
$ws.Range("K73").Value = 438598.74
$ws.Range("M73").Value = -437662.74
# row 74
$ws.Range("H74").Value = 8527
$ws.Range("I74").Value = 6780.8335
$ws.Range("J74").Value = 19004
$ws.Range("K74").Value = 6780.8335
$ws.Range("L74").Value = 19004
$ws.Range("M74").Value = -5844.8335
$ws.Range("N74").Value = -20876
# row 77
$ws.Range("H77").Value = 8527
$ws.Range("I77").Value = 6780.8335
$ws.Range("J77").Value = 19004
$ws.Range("K77").Value = 33904.1675
$ws.Range("L77").Value = 95020
$ws.Range("M77").Value = -29224.1675
$ws.Range("N77").Value = -104380

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 4061.303
$ws.Range("I61").Value = 3393.1785
$ws.Range("J61").Value = 7802.8
$ws.Range("K61").Value = 3393.1785
$ws.Range("L61").Value = 7802.8
$ws.Range("M61").Value = -3181.1785
$ws.Range("N61").Value = -8226.799999999999
# row 74
$ws.Range("H74").Value = 25645812
$ws.Range("I74").Value = 37040440
$ws.Range("J74").Value = 7903.5
$ws.Range("K74").Value = 37040440
$ws.Range("L74").Value = 7903.5
$ws.Range("M74").Value = -37039566
$ws.Range("N74").Value = -9651.5
# row 77
$ws.Range("H77").Value = 25645812
$ws.Range("I77").Value = 37040440
$ws.Range("J77").Value = 7903.5
$ws.Range("K77").Value = 185202200
$ws.Range("L77").Value = 39517.5
$ws.Range("M77").Value = -185197832
$ws.Range("N77").Value = -48253.5
# row 97
$ws.Range("H97").Value = 1004.7778
$ws.Range("J97").Value = 655.6667
$ws.Range("L97").Value = 655.6667
$ws.Range("N97").Value = -1647.6667
# row 132
$ws.Range("H132").Value = 3073.805
$ws.Range("I132").Value = 2208.3845
$ws.Range("J132").Value = 4573.8667
$ws.Range("K132").Value = 6625.1535
$ws.Range("L132").Value = 13721.6001
$ws.Range("M132").Value = -4095.1535
$ws.Range("N132").Value = -18781.6001
# row 136
$ws.Range("H136").Value = 4061.303
$ws.Range("I136").Value = 3393.1785
$ws.Range("J136").Value = 7802.8
$ws.Range("K136").Value = 10179.5355
$ws.Range("L136").Value = 23408.4
$ws.Range("M136").Value = -7629.5355
$ws.Range("N136").Value = -28508.4

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = $null
$ws.Range("N96").Value = $null
# row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = $null
# row 134
$ws.Range("H134").Value = 2016.2
$ws.Range("I134").Value = 1272.4474
$ws.Range("K134").Value = 3817.3422
$ws.Range("M134").Value = -1282.3422
# row 138
$ws.Range("H138").Value = 65382
$ws.Range("J138").Value = 65382
$ws.Range("L138").Value = 65382
$ws.Range("N138").Value = -75662

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 12
$ws.Range("H12").Value = 3001.6667
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 3002.5
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3002.5
$ws.Range("M12").Value = -2830
$ws.Range("N12").Value = -3342.5
# row 16
$ws.Range("H16").Value = 1378.7693
$ws.Range("I16").Value = 697
$ws.Range("J16").Value = 2469.6
$ws.Range("K16").Value = 697
$ws.Range("L16").Value = 2469.6
$ws.Range("M16").Value = -410
$ws.Range("N16").Value = -3043.6
# row 31
$ws.Range("H31").Value = 58582.105
$ws.Range("I31").Value = 2469
$ws.Range("K31").Value = 2469
$ws.Range("M31").Value = -2174
# row 34
$ws.Range("H34").Value = 58582.105
$ws.Range("I34").Value = 2469
$ws.Range("K34").Value = 2469
$ws.Range("M34").Value = -2267
# row 93
$ws.Range("H93").Value = 24999
$ws.Range("I93").Value = 24999
$ws.Range("K93").Value = 24999
$ws.Range("M93").Value = -23127
# row 106
$ws.Range("H106").Value = 52745.668
$ws.Range("J106").Value = 52745.668
$ws.Range("L106").Value = 52745.668
$ws.Range("N106").Value = -55269.668
# row 113
$ws.Range("H113").Value = 1378.7693
$ws.Range("I113").Value = 697
$ws.Range("J113").Value = 2469.6
$ws.Range("K113").Value = 697
$ws.Range("L113").Value = 2469.6
$ws.Range("M113").Value = 1473
$ws.Range("N113").Value = -6809.6
# row 122
$ws.Range("H122").Value = 4758.1055
$ws.Range("I122").Value = 1460.5
$ws.Range("K122").Value = 4381.5
$ws.Range("M122").Value = -1931.5
# row 134
$ws.Range("H134").Value = 2262.7297
$ws.Range("I134").Value = 1414.4828
$ws.Range("J134").Value = 5337.625
$ws.Range("K134").Value = 4243.4484
$ws.Range("L134").Value = 16012.875
$ws.Range("M134").Value = -1708.4484
$ws.Range("N134").Value = -21082.875

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 57038950
$ws.Range("I4").Value = 109000300
$ws.Range("J4").Value = 9801361
$ws.Range("K4").Value = 327000900
$ws.Range("L4").Value = 29404083
$ws.Range("M4").Value = -327000788
$ws.Range("N4").Value = -29404307
# row 88
$ws.Range("H88").Value = 9125.5
$ws.Range("I88").Value = 3235
$ws.Range("J88").Value = 15016
$ws.Range("K88").Value = 9705
$ws.Range("L88").Value = 45048
$ws.Range("M88").Value = -9277
$ws.Range("N88").Value = -45904
# row 91
$ws.Range("H91").Value = 9125.5
$ws.Range("I91").Value = 3235
$ws.Range("J91").Value = 15016
$ws.Range("K91").Value = 9705
$ws.Range("L91").Value = 45048
$ws.Range("M91").Value = -8223
$ws.Range("N91").Value = -48012
# row 98
$ws.Range("H98").Value = 2368.3
$ws.Range("I98").Value = 551
$ws.Range("J98").Value = 2570.2222
$ws.Range("K98").Value = 1653
$ws.Range("L98").Value = 7710.6666
$ws.Range("M98").Value = -155
$ws.Range("N98").Value = -10706.6666
# row 122
$ws.Range("H122").Value = 8026382
$ws.Range("I122").Value = 8547598
$ws.Range("J122").Value = 7939512.5
$ws.Range("K122").Value = 76928382
$ws.Range("L122").Value = 71455612.5
$ws.Range("M122").Value = -76925932
$ws.Range("N122").Value = -71460512.5
# row 129
$ws.Range("H129").Value = 13895406
$ws.Range("J129").Value = 33346940
$ws.Range("L129").Value = 100040820
$ws.Range("N129").Value = -100050820
# row 138
$ws.Range("H138").Value = 1955
$ws.Range("I138").Value = 1476.6666
$ws.Range("K138").Value = 4429.9998
$ws.Range("M138").Value = 710.0002000000004
# row 139
$ws.Range("H139").Value = 4174.2104
$ws.Range("I139").Value = 2980.625
$ws.Range("J139").Value = 5042.273
$ws.Range("K139").Value = 8941.875
$ws.Range("L139").Value = 15126.819
$ws.Range("M139").Value = -3801.875
$ws.Range("N139").Value = -25406.819

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 113
$ws.Range("H113").Value = 2806.0833
$ws.Range("I113").Value = 1909.7333
$ws.Range("K113").Value = 1909.7333
$ws.Range("M113").Value = 260.2666999999999
# row 122
$ws.Range("H122").Value = 9211.087
$ws.Range("I122").Value = 9163.388999999999
$ws.Range("K122").Value = 27490.167
$ws.Range("M122").Value = -25040.167

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 7436.6665
$ws.Range("J46").Value = 8699.299999999999
$ws.Range("L46").Value = 8699.299999999999
$ws.Range("N46").Value = -9075.299999999999
# row 92
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
# row 132
$ws.Range("H132").Value = 6846.737
$ws.Range("I132").Value = 6562.875
$ws.Range("J132").Value = 7053.1816
$ws.Range("K132").Value = 19688.625
$ws.Range("L132").Value = 21159.5448
$ws.Range("M132").Value = -17158.625
$ws.Range("N132").Value = -26219.5448
# row 136
$ws.Range("H136").Value = 3196.2878
$ws.Range("J136").Value = 4327.077
$ws.Range("L136").Value = 12981.231
$ws.Range("N136").Value = -18081.231

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2221.5945
$ws.Range("I132").Value = 1610.5172
$ws.Range("J132").Value = 4436.75
$ws.Range("K132").Value = 4831.5516
$ws.Range("L132").Value = 13310.25
$ws.Range("M132").Value = -2301.5516
$ws.Range("N132").Value = -18370.25
